# Add newly registered players to the roster
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newPlayers = @(
    @{ Name = "Hanj Manfred Elakie Ngalle"; Class = "B1B"; Phone = 650469243 },
    @{ Name = "Kouayep Wanko Ruchi";        Class = "B1B"; Phone = 695353905 },
    @{ Name = "Yvan Longo Gaetan Simon";    Class = "B1A"; Phone = 691674935 },
    @{ Name = "Nzoga Gilbert Boris";        Class = "L1E"; Phone = 657046719 }
)

$startRow = 9
for ($i = 0; $i -lt $newPlayers.Count; $i++) {
    $row = $startRow + $i
    $player = $newPlayers[$i]
    $ws.Cells.Item($row, 1).Value = $player.Name
    $ws.Cells.Item($row, 2).Value = $player.Class
    $ws.Cells.Item($row, 3).Value = $player.Phone
}

$lastRow = $startRow + $newPlayers.Count - 1
$ws.Range("D$lastRow").Select() | Out-Null
